$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New section "漫画目录" (comic chapter-list route) added at the bottom of the
# sheet, rows 65-67, mirroring the look of the other route blocks:
#   row 65 -> section header   (copy formatting from row 61 header block)
#   row 66 -> "接口" route row (copy formatting from row 54 route row, which
#             carries the hyperlink-style formatting used on B2/B25/B36/B54)
#   row 67 -> "参数" param row (copy formatting from row 63 param row, the
#             bordered / thick-bottom closing row of a block)
# ---------------------------------------------------------------------------

$ws.Range("A61:F61").Copy() | Out-Null
$ws.Range("A65:F65").PasteSpecial(-4122) | Out-Null

$ws.Range("A54:F54").Copy() | Out-Null
$ws.Range("A66:F66").PasteSpecial(-4122) | Out-Null

$ws.Range("A63:F63").Copy() | Out-Null
$ws.Range("A67:F67").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# Row heights to match the new block (33pt, custom height)
$ws.Rows.Item(65).RowHeight = 33
$ws.Rows.Item(66).RowHeight = 33
$ws.Rows.Item(67).RowHeight = 33

# Merge the header / route cells
$ws.Range("A65:F65").Merge() | Out-Null
$ws.Range("B66:F66").Merge() | Out-Null

# ---- values ---------------------------------------------------------------
$ws.Range("A65").Value = "漫画目录"

$ws.Range("A66").Value = "接口"
$ws.Range("B66").Value = "https://www.kanman.com/api/getchapterlist"

$ws.Range("A67").Value = "参数"
$ws.Range("B67").Value = "comic_id"
$ws.Range("C67").Value = "是"
$ws.Range("D67").Value = "漫画id"
$ws.Range("E67").Value2 = 106619

# Hyperlink on the new route cell, same pattern as the other route links
$ws.Hyperlinks.Add($ws.Range("B66"), "https://www.kanman.com/api/getchapterlist") | Out-Null

# ---- view / selection state -------------------------------------------
$ws.Range("B66:F66").Select()
